# Auto-generated edit script: append new scrape rows as of 2025-10-29 12:39:32 JST
# Rebuilds sheet1 ("ランサーズ") rows 2..17 with the refreshed scrape data, grows
# column D width 30 -> 32, and rewires the F-column hyperlinks to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop existing hyperlinks up front -- row shifting does not move hyperlink
# anchors along with cell data, so we rebuild them from scratch below.
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(2, 2).Value = '【業務効率化】生成AIを活用したシステム開発依頼'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5422740'
$ws.Cells.Item(2, 7).Value = 445
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆効率化,開発'

# Row 3
$ws.Cells.Item(3, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(3, 2).Value = '【急募】OpenAIを活用した英語力診断ツールの開発依頼'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5423046'
$ws.Cells.Item(3, 7).Value = 408
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◆ツール,開発'

# Row 4
$ws.Cells.Item(4, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(4, 2).Value = '【急募】業務効率化・生成AI実装のAIエンジニアパートナー募集'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5422760'
$ws.Cells.Item(4, 7).Value = 385
$ws.Cells.Item(4, 8).Value = '🔥AI,Ai ◆効率化'

# Row 5
$ws.Cells.Item(5, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(5, 2).Value = '【機密性の高いノウハウを含む】サーバーレスAI分析システム構築(MVP開発と拡張性確保)'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5422386'
$ws.Cells.Item(5, 7).Value = 383
$ws.Cells.Item(5, 8).Value = '🔥AI,Ai ◆開発'

# Row 6
$ws.Cells.Item(6, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(6, 2).Value = '【自動化】EAを証券口座・VPSに接続する作業の効率化依頼'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5422331'
$ws.Cells.Item(6, 7).Value = 145
$ws.Cells.Item(6, 8).Value = '◆効率化,自動化'

# Row 7
$ws.Cells.Item(7, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(7, 2).Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Cells.Item(7, 7).Value = 135
$ws.Cells.Item(7, 8).Value = '◆ツール,スクレイピング ◇サイト'

# Row 8
$ws.Cells.Item(8, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(8, 2).Value = '【動画解析】動作比較アルゴリズム開発者を募集します'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5422314'
$ws.Cells.Item(8, 7).Value = 68
$ws.Cells.Item(8, 8).Value = '◆開発'

# Row 9
$ws.Cells.Item(9, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(9, 2).Value = 'アマゾンの返品レポートより返品理由のポップアップ文字までダウンロードしてエクセルに書き出すツール'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5422652'
$ws.Cells.Item(9, 7).Value = 65
$ws.Cells.Item(9, 8).Value = '◆ツール'

# Row 10
$ws.Cells.Item(10, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(10, 2).Value = '【急募】Accessシステム改修・CSV読込・MySQLクラウド化【出張希望】'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5422936'
$ws.Cells.Item(10, 7).Value = 48
$ws.Cells.Item(10, 8).Value = '◇MySQL'

# Row 11
$ws.Cells.Item(11, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(11, 2).Value = '【急募】WordPressにe-SCOTT決済機能を導入'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5422908'
$ws.Cells.Item(11, 7).Value = 33
$ws.Cells.Item(11, 8).Value = '○WordPress'

# Row 12
$ws.Cells.Item(12, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(12, 2).Value = '【高額成功報酬】レガシー基幹システムのバイナリ解析とパッチ作成'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5415980'
$ws.Cells.Item(12, 7).Value = 40
$ws.Cells.Item(12, 8).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(13, 2).Value = '【急募】東京でのマクロ構築依頼!スキルを活かしませんか?'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5423114'
$ws.Cells.Item(13, 7).Value = 25
$ws.Cells.Item(13, 8).ClearContents()

# Row 14
$ws.Cells.Item(14, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(14, 2).Value = '【プロトタイプ基板】重量測定用基板の設計・製造依頼'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5422916'
$ws.Cells.Item(14, 7).Value = 18
$ws.Cells.Item(14, 8).ClearContents()

# Row 15
$ws.Cells.Item(15, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(15, 2).Value = '【オンライン講師募集】HTML・CSSの基礎を実践的に教えていただける方'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5422660'
$ws.Cells.Item(15, 7).Value = 18
$ws.Cells.Item(15, 8).ClearContents()

# Row 16
$ws.Cells.Item(16, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(16, 2).Value = '【急募】YouTubeの音楽配信構築の依頼です'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5420233'
$ws.Cells.Item(16, 7).Value = 13
$ws.Cells.Item(16, 8).ClearContents()

# Row 17
$ws.Cells.Item(17, 1).Value = '2025-10-29 12:39:32'
$ws.Cells.Item(17, 2).Value = '限定公開 限定公開の仕事'
$ws.Cells.Item(17, 3).Value = 'システム開発'
$ws.Cells.Item(17, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(17, 5).Value = '期限情報なし'
$ws.Cells.Item(17, 6).Value = 'https://www.lancers.jp/work/detail/5421782'
$ws.Cells.Item(17, 7).Value = 10
$ws.Cells.Item(17, 8).ClearContents()

# Re-add the F-column hyperlinks (one per data row) and restore the shared
# built-in "Hyperlink" cell style, reusing the style slot already on the sheet.
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5422740')
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5423046')
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5422760')
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5422386')
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5422331')
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5251319')
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5422314')
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5422652')
$ws.Cells.Item(9, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5422936')
$ws.Cells.Item(10, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5422908')
$ws.Cells.Item(11, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5415980')
$ws.Cells.Item(12, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5423114')
$ws.Cells.Item(13, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5422916')
$ws.Cells.Item(14, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5422660')
$ws.Cells.Item(15, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), 'https://www.lancers.jp/work/detail/5420233')
$ws.Cells.Item(16, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), 'https://www.lancers.jp/work/detail/5421782')
$ws.Cells.Item(17, 6).Style = "Hyperlink"

# Widen column D (price) from 30 to 32 characters; ColumnWidth round-trips
# through Excel's pixel-based storage with a constant +5/6 offset, so back
# the input off by that amount to land exactly on a stored width of 32.
$ws.Columns.Item(4).ColumnWidth = 31.166666666666668

$ws.Range("A1").Select()

